$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.49950888782199
$ws.Range("C2").Value = 4.887934018899122
$ws.Range("E2").Value = 10.84314820719703
$ws.Range("F2").Value = 48.03215854380844
$ws.Range("G2").Value = 3.756706411964176
$ws.Range("I2").Value = 35.10461217398687
$ws.Range("J2").Value = 10.25428702937504
$ws.Range("K2").Value = 15.40012242737566
$ws.Range("L2").Value = 11.55002993862346
$ws.Range("M2").Value = 17.47034419889365
$ws.Range("N2").Value = 24.58050386192985
$ws.Range("B3").Value = 17.37953526142808
$ws.Range("C3").Value = 4.732791445289155
$ws.Range("E3").Value = 10.85391462275479
$ws.Range("F3").Value = 48.02260685275188
$ws.Range("G3").Value = 3.759423490622539
$ws.Range("I3").Value = 35.15757813093179
$ws.Range("J3").Value = 10.26418018169992
$ws.Range("K3").Value = 15.32004611744008
$ws.Range("L3").Value = 11.55903775805177
$ws.Range("M3").Value = 17.46502518455279
$ws.Range("N3").Value = 24.63573892978388
$ws.Range("B4").Value = 17.30963935431379
$ws.Range("C4").Value = 4.636536189421834
$ws.Range("E4").Value = 10.86139474191786
$ws.Range("F4").Value = 48.02574198225521
$ws.Range("G4").Value = 3.761179866074676
$ws.Range("I4").Value = 35.19540732490741
$ws.Range("J4").Value = 10.2704920358906
$ws.Range("K4").Value = 15.27421455790508
$ws.Range("L4").Value = 11.56605586058928
$ws.Range("M4").Value = 17.46460891160368
$ws.Range("N4").Value = 24.67157768027342
$ws.Range("B5").Value = 17.28213039872448
$ws.Range("C5").Value = 4.597130245005647
$ws.Range("E5").Value = 10.86466218338205
$ws.Range("F5").Value = 48.0292849208512
$ws.Range("G5").Value = 3.761917825765259
$ws.Range("I5").Value = 35.21215651781753
$ws.Range("J5").Value = 10.27312411037364
$ws.Range("K5").Value = 15.25639205977417
$ws.Range("L5").Value = 11.56929046473482
$ws.Range("M5").Value = 17.46515799838344
$ws.Range("N5").Value = 24.68666684590406
$ws.Range("B6").Value = 17.27762214097799
$ws.Range("C6").Value = 4.59057809647588
$ws.Range("E6").Value = 10.8652179977866
$ws.Range("F6").Value = 48.03001008872303
$ws.Range("G6").Value = 3.762041707779475
$ws.Range("I6").Value = 35.21501820881458
$ws.Range("J6").Value = 10.2735647918115
$ws.Range("K6").Value = 15.25348466854075
$ws.Range("L6").Value = 11.56985021624132
$ws.Range("M6").Value = 17.46529263996369
$ws.Range("N6").Value = 24.68920167967166
$ws.Range("B7").Value = 17.30926438056162
$ws.Range("C7").Value = 4.636005385346595
$ws.Range("E7").Value = 10.86143791932112
$ws.Range("F7").Value = 48.02578058916709
$ws.Range("G7").Value = 3.761189728380726
$ws.Range("I7").Value = 35.19562781293752
$ws.Range("J7").Value = 10.27052728994995
$ws.Range("K7").Value = 15.27397071921785
$ws.Range("L7").Value = 11.56609796570321
$ws.Range("M7").Value = 17.46461340408512
$ws.Range("N7").Value = 24.67177921503455
$ws.Range("B8").Value = 17.45737631453924
$ws.Range("C8").Value = 4.834688154529508
$ws.Range("E8").Value = 10.8466803589343
$ws.Range("F8").Value = 48.02699884931592
$ws.Range("G8").Value = 3.757625024956567
$ws.Range("I8").Value = 35.12177249834693
$ws.Range("J8").Value = 10.25764905866651
$ws.Range("K8").Value = 15.37182946085961
$ws.Range("L8").Value = 11.55282758134634
$ws.Range("M8").Value = 17.46792014574847
$ws.Range("N8").Value = 24.59914998564162
$ws.Range("B9").Value = 17.77647879015111
$ws.Range("C9").Value = 5.213572081374961
$ws.Range("E9").Value = 10.8246140567587
$ws.Range("F9").Value = 48.1006591666706
$ws.Range("G9").Value = 3.751330116676649
$ws.Range("I9").Value = 35.01910852905304
$ws.Range("J9").Value = 10.23426736937343
$ws.Range("K9").Value = 15.58946070313552
$ws.Range("L9").Value = 11.53857348734287
$ws.Range("M9").Value = 17.49690934655285
$ws.Range("N9").Value = 24.47195834763625
$ws.Range("B10").Value = 18.02662734416986
$ws.Range("C10").Value = 5.481967954635303
$ws.Range("E10").Value = 10.81255745055873
$ws.Range("F10").Value = 48.19798063927271
$ws.Range("G10").Value = 3.74712446856062
$ws.Range("I10").Value = 34.9694494877444
$ws.Range("J10").Value = 10.21821408062382
$ws.Range("K10").Value = 15.76397843319624
$ws.Range("L10").Value = 11.53523081208338
$ws.Range("M10").Value = 17.53176925213479
$ws.Range("N10").Value = 24.38775037198593
$ws.Range("B11").Value = 18.14343629564108
$ws.Range("C11").Value = 5.601280520382866
$ws.Range("E11").Value = 10.80796742957017
$ws.Range("F11").Value = 48.25156004683161
$ws.Range("G11").Value = 3.745301228287851
$ws.Range("I11").Value = 34.95246324566297
$ws.Range("J11").Value = 10.21115184952372
$ws.Range("K11").Value = 15.84629938912607
$ws.Range("L11").Value = 11.53524738305041
$ws.Range("M11").Value = 17.55053261267133
$ws.Range("N11").Value = 24.35143845548361
$ws.Range("B12").Value = 18.18806744709924
$ws.Range("C12").Value = 5.64601592930628
$ws.Range("E12").Value = 10.80635730616386
$ws.Range("F12").Value = 48.27317901520944
$ws.Range("G12").Value = 3.744623669471139
$ws.Range("I12").Value = 34.94683719222569
$ws.Range("J12").Value = 10.20851188674683
$ws.Range("K12").Value = 15.87787169186921
$ws.Range("L12").Value = 11.53547363402824
$ws.Range("M12").Value = 17.55805163453269
$ws.Range("N12").Value = 24.33797423673607
$ws.Range("B13").Value = 18.17843819337314
$ws.Range("C13").Value = 5.636401855185627
$ws.Range("E13").Value = 10.80669839144508
$ws.Range("F13").Value = 48.26846400124936
$ws.Range("G13").Value = 3.744769022909026
$ws.Range("I13").Value = 34.94801300004482
$ws.Range("J13").Value = 10.20907892587666
$ws.Range("K13").Value = 15.87105464478278
$ws.Range("L13").Value = 11.53541514115636
$ws.Range("M13").Value = 17.55641394516253
$ws.Range("N13").Value = 24.34086127126479
$ws.Range("B14").Value = 18.14710035318439
$ws.Range("C14").Value = 5.604970105491157
$ws.Range("E14").Value = 10.80783240168043
$ws.Range("F14").Value = 48.25331205953619
$ws.Range("G14").Value = 3.745245227722143
$ws.Range("I14").Value = 34.95198422638847
$ws.Range("J14").Value = 10.21093397107795
$ws.Range("K14").Value = 15.84888899325432
$ws.Range("L14").Value = 11.53526159459496
$ws.Range("M14").Value = 17.55114293941488
$ws.Range("N14").Value = 24.35032501082766
$ws.Range("B15").Value = 18.12795580828251
$ws.Range("C15").Value = 5.58565794019937
$ws.Range("E15").Value = 10.80854366875329
$ws.Range("F15").Value = 48.24420393125286
$ws.Range("G15").Value = 3.745538589980953
$ws.Range("I15").Value = 34.95452172600455
$ws.Range("J15").Value = 10.2120747064677
$ws.Range("K15").Value = 15.83536319293909
$ws.Range("L15").Value = 11.53519615726668
$ws.Range("M15").Value = 17.54796805222718
$ws.Range("N15").Value = 24.3561590961994
$ws.Range("B16").Value = 18.01905126203211
$ws.Range("C16").Value = 5.474110606010878
$ws.Range("E16").Value = 10.81287536762612
$ws.Range("F16").Value = 48.19466566389645
$ws.Range("G16").Value = 3.747245425303309
$ws.Range("I16").Value = 34.97067238572231
$ws.Range("J16").Value = 10.2186804335244
$ws.Range("K16").Value = 15.75865565957803
$ws.Range("L16").Value = 11.53526057144332
$ws.Range("M16").Value = 17.53060114787452
$ws.Range("N16").Value = 24.39016353486221
$ws.Range("B17").Value = 17.95298909454678
$ws.Range("C17").Value = 5.404933378768228
$ws.Range("E17").Value = 10.81576141026582
$ws.Range("F17").Value = 48.16665351635067
$ws.Range("G17").Value = 3.748315497411123
$ws.Range("I17").Value = 34.98201596713844
$ws.Range("J17").Value = 10.22279426367312
$ws.Range("K17").Value = 15.71233326699314
$ws.Range("L17").Value = 11.53569305859637
$ws.Range("M17").Value = 17.52068848157394
$ws.Range("N17").Value = 24.41153467372566
$ws.Range("B18").Value = 17.915277848401
$ws.Range("C18").Value = 5.364885018123351
$ws.Range("E18").Value = 10.81750562658763
$ws.Range("F18").Value = 48.15141832588949
$ws.Range("G18").Value = 3.74893944333684
$ws.Range("I18").Value = 34.98906794938315
$ws.Range("J18").Value = 10.22518307982807
$ws.Range("K18").Value = 15.68596721712007
$ws.Range("L18").Value = 11.5360865682496
$ws.Range("M18").Value = 17.51526071242585
$ws.Range("N18").Value = 24.42401456934599
$ws.Range("B19").Value = 17.90255965221792
$ws.Range("C19").Value = 5.351282240014036
$ws.Range("E19").Value = 10.81811067553318
$ws.Range("F19").Value = 48.14641076676921
$ws.Range("G19").Value = 3.749152157338552
$ws.Range("I19").Value = 34.99154620332394
$ws.Range("J19").Value = 10.22599578913807
$ws.Range("K19").Value = 15.67708843603444
$ws.Range("L19").Value = 11.53624469684519
$ws.Range("M19").Value = 17.51347009532986
$ws.Range("N19").Value = 24.42827231551045
$ws.Range("B20").Value = 17.95999218336619
$ws.Range("C20").Value = 5.412324610465543
$ws.Range("E20").Value = 10.81544547197591
$ws.Range("F20").Value = 48.16954478644685
$ws.Range("G20").Value = 3.748200710426315
$ws.Range("I20").Value = 34.9807538300814
$ws.Range("J20").Value = 10.22235399692484
$ws.Range("K20").Value = 15.71723581896559
$ws.Range("L20").Value = 11.5356320455478
$ws.Range("M20").Value = 17.52171539965958
$ws.Range("N20").Value = 24.40924024860204
$ws.Range("B21").Value = 18.1562945165498
$ws.Range("C21").Value = 5.614214807361593
$ws.Range("E21").Value = 10.80749584603047
$ws.Range("F21").Value = 48.25772654274857
$ws.Range("G21").Value = 3.74510500628354
$ws.Range("I21").Value = 34.95079589674969
$ws.Range("J21").Value = 10.21038816890294
$ws.Range("K21").Value = 15.85538893233728
$ws.Range("L21").Value = 11.5353007337457
$ws.Range("M21").Value = 17.55267996654403
$ws.Range("N21").Value = 24.34753751464121
$ws.Range("B22").Value = 18.28689326538007
$ws.Range("C22").Value = 5.743546193792012
$ws.Range("E22").Value = 10.80304619638412
$ws.Range("F22").Value = 48.32310349104878
$ws.Range("G22").Value = 3.743156726621438
$ws.Range("I22").Value = 34.93591623976944
$ws.Range("J22").Value = 10.20276792482404
$ws.Range("K22").Value = 15.94799559667021
$ws.Range("L22").Value = 11.53636584101791
$ws.Range("M22").Value = 17.57532661543601
$ws.Range("N22").Value = 24.30887985103937
$ws.Range("B23").Value = 18.21699130726052
$ws.Range("C23").Value = 5.674772759344414
$ws.Range("E23").Value = 10.80535301744557
$ws.Range("F23").Value = 48.28750511594507
$ws.Range("G23").Value = 3.744189725670083
$ws.Range("I23").Value = 34.94342769675421
$ws.Range("J23").Value = 10.20681675814132
$ws.Range("K23").Value = 15.89836539729654
$ws.Range("L23").Value = 11.53568047561074
$ws.Range("M23").Value = 17.56302062300124
$ws.Range("N23").Value = 24.32935966043259
$ws.Range("B24").Value = 17.95682524699116
$ws.Range("C24").Value = 5.408983895306005
$ws.Range("E24").Value = 10.81558804291549
$ws.Range("F24").Value = 48.16823493448197
$ws.Range("G24").Value = 3.74825257836905
$ws.Range("I24").Value = 34.9813227900072
$ws.Range("J24").Value = 10.22255296760075
$ws.Range("K24").Value = 15.71501854707403
$ws.Range("L24").Value = 11.53565917822766
$ws.Range("M24").Value = 17.52125028505651
$ws.Range("N24").Value = 24.41027695584614
$ws.Range("B25").Value = 17.68726261070109
$ws.Range("C25").Value = 5.112596534952814
$ws.Range("E25").Value = 10.82985128888759
$ws.Range("F25").Value = 48.07312274435973
$ws.Range("G25").Value = 3.752959095000199
$ws.Range("I25").Value = 35.04236041520695
$ws.Range("J25").Value = 10.24039395877489
$ws.Range("K25").Value = 15.52793649874665
$ws.Range("L25").Value = 11.54117399016338
$ws.Range("M25").Value = 17.48667274015308
$ws.Range("N25").Value = 24.50474110985745
